# Restore revision: update the "From" (min hour) value for rule R30 in the
# Rules table from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
